$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "51.735.10"
$ws.Range("E2").Value = "  +2.01%  "

# Row 3
$ws.Range("D3").Value = "3.039.77"
$ws.Range("E3").Value = "  +4.13%  "

# Row 4
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
Set-TextValue "D5" "381.11"
$ws.Range("E5").Value = "  +1.69%  "

# Row 6
Set-TextValue "D6" "103.63"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7
Set-TextValue "D7" "0.546"
$ws.Range("E7").Value = "  +2.25%  "

# Row 8
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
Set-TextValue "D9" "0.598"
$ws.Range("E9").Value = "  +4.00%  "

# Row 10
Set-TextValue "D10" "36.96"
$ws.Range("E10").Value = "  +3.91%  "

# Row 11
$ws.Range("E11").Value = "  -0.05%  "

# Row 12
Set-TextValue "D12" "0.0861"
$ws.Range("E12").Value = "  +1.86%  "

# Row 13
$ws.Range("D13").Value = "3.495.96"
$ws.Range("E13").Value = "  +3.31%  "

# Row 14
Set-TextValue "D14" "18.63"
$ws.Range("E14").Value = "  +3.84%  "

# Row 15
Set-TextValue "D15" "7.78"
$ws.Range("E15").Value = "  +2.41%  "

# Row 16
$ws.Range("D16").Value = "3.044.83"
$ws.Range("E16").Value = "  +3.26%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "1.00"
$ws.Range("E17").Value = "  +1.12%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D18" "10.89"
$ws.Range("E18").Value = "  -9.57%  "

# Row 19
$ws.Range("D19").Value = "51.832.57"
$ws.Range("E19").Value = "  +2.25%  "

# Row 20
$ws.Range("E20").Value = "  +3.01%  "

# Row 21
Set-TextValue "D21" "12.58"
$ws.Range("E21").Value = "  +2.82%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  +2.72%  "

# Row 23
Set-TextValue "D23" "70.44"
$ws.Range("E23").Value = "  +1.79%  "

# Row 24
Set-TextValue "D24" "269.04"
$ws.Range("E24").Value = "  +1.34%  "

# Row 25
$ws.Range("E25").Value = "  +1.88%  "

# Row 26
$ws.Range("E26").Value = "  +6.39%  "

# Row 27
Set-TextValue "D27" "7.57"
$ws.Range("E27").Value = "  +6.47%  "

# Row 28
Set-TextValue "D28" "0.171"
$ws.Range("E28").Value = "  +6.25%  "

# Row 29
Set-TextValue "D29" "26.28"
$ws.Range("E29").Value = "  +3.76%  "

# Row 30
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
$ws.Range("E31").Value = "  +1.87%  "

# Row 32
Set-TextValue "D32" "10.36"
$ws.Range("E32").Value = "  +4.91%  "

# Row 33
Set-TextValue "D33" "34.55"
$ws.Range("E33").Value = "  +4.47%  "

# Row 34
Set-TextValue "D34" "51.22"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D35" "0.0457"
$ws.Range("E35").Value = "  +6.70%  "

# Row 36
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D36" "2.06"
$ws.Range("E36").Value = "  +0.79%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
Set-TextValue "D38" "3.33"
$ws.Range("E38").Value = "  +9.20%  "

# Row 39
Set-TextValue "D39" "17.17"
$ws.Range("E39").Value = "  +5.61%  "

# Row 40
Set-TextValue "D40" "2.62"
$ws.Range("E40").Value = "  +8.93%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D41" "0.285"
$ws.Range("E41").Value = "  +10.79%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "1.87"
$ws.Range("E42").Value = "  +4.81%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D43" "0.117"
$ws.Range("E43").Value = "  +1.78%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "126.09"
$ws.Range("E44").Value = "  +5.49%  "

# Row 45
Set-TextValue "D45" "3.80"
$ws.Range("E45").Value = "  +13.89%  "

# Row 46
Set-TextValue "D46" "21.92"
$ws.Range("E46").Value = "  +5.41%  "

# Row 47
Set-TextValue "D47" "2.04"
$ws.Range("E47").Value = "  +0.75%  "

# Row 48
Set-TextValue "D48" "2.39"
$ws.Range("E48").Value = "  +2.45%  "

# Row 49
$ws.Range("D49").Value = "2.047.06"
$ws.Range("E49").Value = "  +3.06%  "

# Row 50
$ws.Range("D50").Value = "3.336.04"
$ws.Range("E50").Value = "  +3.96%  "

# Row 51
Set-TextValue "D51" "0.0327"
$ws.Range("E51").Value = "  +4.30%  "

"done"